# Update MathSAT data set: shift the "Year" column (A2:A40) on the "Data"
# sheet forward by 19 years (1967-2005 -> 1986-2024), bring the formatting
# of the newly-shifted rows (A32:A40) in line with the rest of the column,
# and make "Data" the active sheet with the selection left on the first
# empty row below the table (A41).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

for ($r = 2; $r -le 40; $r++) {
    $cell = $dataSheet.Cells.Item($r, 1)
    $year = $cell.Value()
    $cell.Value = $year + 19
}

# A2:A31 already carry the bordered "year" style; copy it onto A32:A40 so
# the whole column looks consistent after the shift.
$dataSheet.Range("A31").Copy()
$dataSheet.Range("A32:A40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Make "Data" the active/selected sheet, with the selection sitting just
# below the table.
$dataSheet.Activate()
$dataSheet.Range("A41").Select()
